# Auto commit at 2025-10-10 9:45:54.25
# Updates the Metrics sheet source values (which ripple via formulas into the
# "today" sheet), and moves the active-tab/selection state from Metrics to
# the "today" sheet.

$wb = $excel.ActiveWorkbook

$wsMetrics = $wb.Worksheets.Item("Metrics")
$wsToday   = $wb.Worksheets.Item("today")

# --- Update the source metric values on the Metrics sheet ---------------
$wsMetrics.Range("B2").Value  = 125373.59000000001
$wsMetrics.Range("B3").Value  = 105712.89000000001
$wsMetrics.Range("B4").Value  = 47073
$wsMetrics.Range("B5").Value  = 4857
$wsMetrics.Range("B6").Value  = 4492505.0600000005
$wsMetrics.Range("B7").Value  = 3795531.5599999996
$wsMetrics.Range("B8").Value  = 1317675.1400000001
$wsMetrics.Range("B9").Value  = 173858
$wsMetrics.Range("B10").Value = 32957828.860999826
$wsMetrics.Range("B11").Value = 31070753.080000006
$wsMetrics.Range("B12").Value = 11599384.029999999
$wsMetrics.Range("B13").Value = 1271485

# --- Move the selection / active-tab state -------------------------------
# Select Metrics!D8 first (leaves that sheet's cached selection at D8,
# without leaving it flagged as the active tab)...
$wsMetrics.Range("D8").Select()

# ...then activate "today" and select B6 last, so it becomes the active
# sheet (activeTab) with tabSelected="1" and the new selection anchored at
# B6.
$wsToday.Activate()
$wsToday.Range("B6").Select()
